# BancoDados.xlsx: "Cadastro" sheet gets a new "UserNameFalha" column (B),
# header renames/capitalization, an updated userName test value, a couple of
# new blank styled placeholder cells, and an updated selection. The
# "Pesquisa" sheet is untouched content-wise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cadastro")

# 1. Insert a new column B ("UserNameFalha"); this shifts former B:K to C:L
#    and correctly carries over widths/styles/blank placeholder cells
#    (A3/K3 -> A3,B3 / L3) automatically.
$ws.Columns.Item(2).Insert()

# 2. Re-apply header row (row 1) with the renamed / capitalized headers.
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "UserNameFalha"
$ws.Range("C1").Value = "UserPass"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "FirstName"
$ws.Range("F1").Value = "LastName"
$ws.Range("G1").Value = "PhoneNumber"
$ws.Range("H1").Value = "City"
$ws.Range("I1").Value = "Address"
$ws.Range("J1").Value = "State"
$ws.Range("K1").Value = "PostalCode"
$ws.Range("L1").Value = "Country"

# 3. Re-apply the data row (row 2) with the new userName value and the new
#    UserNameFalha value; the rest keep their existing text.
$ws.Range("A2").Value = "lucasVi"
$ws.Range("B2").Value = "lucasViadoadasdda"

# 4. Fix up the hyperlinks: column insert does not shift hyperlink ranges
#    in this engine, so remove + recreate them on the now-correct cells,
#    then restore the original (non-hyperlink) cell formatting.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Trocar@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:abdiel.cordeiro@rsinet.com.br") | Out-Null
$ws.Range("A2").Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. New blank styled placeholder cells (D5, K7) mirroring the look of the
#    existing A3/L3 placeholders.
$ws.Range("A3").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("L3").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 6. Column width tweaks for the new/resized columns (closest achievable
#    values given this engine's column-width rounding).
$ws.Columns.Item(2).ColumnWidth = 15.83
$ws.Columns.Item(11).ColumnWidth = 9.6

# 7. Selection moves to E15.
$ws.Range("E15").Select()
